$d = $word.ActiveDocument

# Locate the paragraph that starts the footer block to be removed
# ("Ver no Jupiter Salvar em pdf Salvar em docx") and the paragraph that
# ends it ("... Powered by Jekyll and Github pages. ...").
$startIdx = $null
$endIdx = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "*Ver no Jupiter*") {
        $startIdx = $i
    }
    if ($t -like "*Powered by Jekyll*") {
        $endIdx = $i
    }
}

if ($startIdx -ne $null -and $endIdx -ne $null) {
    # Also remove the blank paragraph that immediately precedes the block,
    # matching the surrounding blank-paragraph pattern used elsewhere in
    # the document.
    $precedingIdx = $startIdx - 1
    if ($precedingIdx -ge 1) {
        $preceding = $d.Paragraphs.Item($precedingIdx)
        if ($preceding.Range.Text.Trim() -eq "") {
            $startIdx = $precedingIdx
        }
    }

    $rangeStart = $d.Paragraphs.Item($startIdx).Range.Start
    $rangeEnd = $d.Paragraphs.Item($endIdx).Range.End
    $r = $d.Range($rangeStart, $rangeEnd)
    $r.Delete()
}
